$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target state for rows 2-19 (Player, Position, Team)
$data = @(
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Ja Morant", "PG", "Memphis Grizzlies"),
    @("Scottie Barnes", "SG,SF,PF", "Toronto Raptors"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Vasilije Micic", "PG,SG", "Charlotte Hornets"),
    @("Bobby Portis", "PF,C", "Milwaukee Bucks"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Isaiah Stewart", "PF,C", "Detroit Pistons"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row = $row + 1
}
